# Update "Förändrad" date column (C2:C7) from serial 45188 (2023-09-19)
# to serial 45189 (2023-09-20), i.e. advance by one day.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

foreach ($row in 2..7) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value2 = $cell.Value2 + 1
}
